# Daily refresh of the cryptos list (prices / 1h volume deltas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new "Price" values look like plain numbers (e.g. "99.97").
# The source data always stores that column as text (so values such as
# "42.909.56" - which use '.' as a thousands separator - round-trip
# correctly), so force a text number format on those specific cells before
# assigning the value. This keeps Excel from re-interpreting them as
# floating point numbers and losing formatting (trailing zeros, etc.).
$textPriceCells = @(
    "D5","D6","D9","D10","D12","D15","D17","D19","D20","D22","D23",
    "D26","D29","D30","D31","D33","D37","D41","D42","D45","D50","D51"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.909.56"
$ws.Range("E2").Value = "  +0.27%  "

# Row 3 - Ethereum
$ws.Range("E3").Value = "  +1.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "315.44"
$ws.Range("E5").Value = "  +0.28%  "

# Row 6 - Solana
$ws.Range("D6").Value = "99.97"
$ws.Range("E6").Value = "  +4.34%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.26%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.04%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  +0.57%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "36.38"
$ws.Range("E10").Value = "  +0.55%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.50%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "7.52"
$ws.Range("E12").Value = "  -0.46%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.971.03"
$ws.Range("E13").Value = "  +1.89%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.65%  "

# Row 15 - now Chainlink (was WrappedEther)
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "15.80"
$ws.Range("E15").Value = "  +3.27%  "

# Row 16 - now WrappedEther (was Chainlink)
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.560.52"
$ws.Range("E16").Value = "  +2.22%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  -0.93%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "42.976.33"
$ws.Range("E18").Value = "  +0.26%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "6.84"
$ws.Range("E19").Value = "  +1.20%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  -1.80%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("E21").Value = "  +0.58%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "69.50"
$ws.Range("E22").Value = "  -0.41%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "249.99"
$ws.Range("E23").Value = "  -1.20%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  +0.60%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -0.36%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "27.15"
$ws.Range("E26").Value = "  +1.81%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.01%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  -1.38%  "

# Row 29 - InjectiveProtocol
$ws.Range("D29").Value = "40.68"
$ws.Range("E29").Value = "  -0.42%  "

# Row 30 - Cosmos
$ws.Range("D30").Value = "10.32"
$ws.Range("E30").Value = "  -0.68%  "

# Row 31 - Monero
$ws.Range("D31").Value = "158.02"
$ws.Range("E31").Value = "  +0.39%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -1.79%  "

# Row 33 - LidoDAOToken
$ws.Range("D33").Value = "3.43"
$ws.Range("E33").Value = "  +3.51%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  -1.51%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  +3.07%  "

# Row 36 - WEMIXToken
$ws.Range("E36").Value = "  -0.36%  "

# Row 37 - Celestia
$ws.Range("D37").Value = "18.88"
$ws.Range("E37").Value = "  -2.83%  "

# Row 38 - ApeXProtocol
$ws.Range("E38").Value = "  +9.05%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +1.26%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +0.24%  "

# Row 41 - EnergySwap
$ws.Range("D41").Value = "23.45"
$ws.Range("E41").Value = "  -0.24%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "4.12"
$ws.Range("E42").Value = "  +8.51%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  -0.41%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  -0.07%  "

# Row 45 - NEARProtocol
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  -2.01%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.008.27"
$ws.Range("E46").Value = "  -2.21%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  +0.12%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "2.821.62"
$ws.Range("E48").Value = "  +1.98%  "

# Row 49 - Algorand
$ws.Range("E49").Value = "  +2.61%  "

# Row 50 - ordi
$ws.Range("D50").Value = "75.13"
$ws.Range("E50").Value = "  -0.57%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "81.82"
$ws.Range("E51").Value = "  -4.10%  "
